$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("LoginPageData")
$ws2 = $wb.Worksheets.Item("DashboardPageData")
$ws3 = $wb.Worksheets.Item("InsuredPageData")

$xlPasteFormats = -4122

# ---------------------------------------------------------------
# 1. Add the new "testSubmissionClearancesFunctionality" block
#    (rows 26-29) to the InsuredPageData sheet, mirroring the
#    layout/styling of the existing "testCheckDuplicateSubmission"
#    block (rows 21-23).
# ---------------------------------------------------------------

# Row 26: section title row -> copy formatting from row 21 (A:B)
$ws3.Range("A21:B21").Copy()
$ws3.Range("A26").PasteSpecial($xlPasteFormats)
$ws3.Range("A26").Value = "testSubmissionClearancesFunctionality"

# Row 27: header row -> copy formatting from row 22 (A:I), plus one
# extra column (J) copied from the same style.
$ws3.Range("A22:I22").Copy()
$ws3.Range("A27").PasteSpecial($xlPasteFormats)
$ws3.Range("I22").Copy()
$ws3.Range("J27").PasteSpecial($xlPasteFormats)

$ws3.Range("A27").Value = "runMode"
$ws3.Range("B27").Value = "product"
$ws3.Range("C27").Value = "applicantName"
$ws3.Range("D27").Value = "website"
$ws3.Range("E27").Value = "email"
$ws3.Range("F27").Value = "brokerId"
$ws3.Range("G27").Value = "agentId"
$ws3.Range("H27").Value = "agencyOfficeId"
$ws3.Range("I27").Value = "functionality"
$ws3.Range("J27").Value = "clearanceText"

# Row 28 & 29: data rows -> copy formatting column by column from the
# closest-matching existing cells so that borders / alignment match.
foreach ($r in @(28, 29)) {
    $ws3.Range("A23").Copy()
    $ws3.Range("A$r").PasteSpecial($xlPasteFormats)

    $ws3.Range("B22").Copy()
    $ws3.Range("B$r").PasteSpecial($xlPasteFormats)

    $ws3.Range("C23").Copy()
    $ws3.Range("C$r").PasteSpecial($xlPasteFormats)

    $ws3.Range("C23").Copy()
    $ws3.Range("D$r").PasteSpecial($xlPasteFormats)

    $ws3.Range("E23").Copy()
    $ws3.Range("E$r").PasteSpecial($xlPasteFormats)

    $ws3.Range("F23").Copy()
    $ws3.Range("F$r").PasteSpecial($xlPasteFormats)

    $ws3.Range("G23").Copy()
    $ws3.Range("G$r").PasteSpecial($xlPasteFormats)

    $ws3.Range("H23").Copy()
    $ws3.Range("H$r").PasteSpecial($xlPasteFormats)

    $ws3.Range("I23").Copy()
    $ws3.Range("I$r").PasteSpecial($xlPasteFormats)

    $ws3.Range("I22").Copy()
    $ws3.Range("J$r").PasteSpecial($xlPasteFormats)
}

$ws3.Range("A28").Value = "Y"
$ws3.Range("B28").Value = "QA Program 5203"
$ws3.Range("C28").Value = "ui"
$ws3.Range("D28").Value = "www.ui.com"
$ws3.Range("E28").Value = "cfessler@profrisk.com"
$ws3.Range("F28").Value = 20217
$ws3.Range("G28").Value = 237
$ws3.Range("H28").Value = 8006
$ws3.Range("I28").Value = "submit"
$ws3.Range("J28").Value = "Test purpose"

$ws3.Range("A29").Value = "N"
$ws3.Range("B29").Value = "QA Program 5203"
$ws3.Range("C29").Value = "ui"
$ws3.Range("D29").Value = "www.ui.com"
$ws3.Range("E29").Value = "cfessler@profrisk.com"
$ws3.Range("F29").Value = 20217
$ws3.Range("G29").Value = 237
$ws3.Range("H29").Value = 8006
$ws3.Range("I29").Value = "cancel"
$ws3.Range("J29").Value = "Test purpose"

# ---------------------------------------------------------------
# 2. Switch the active tab from DashboardPageData to
#    InsuredPageData and update the selected cell there.
# ---------------------------------------------------------------
$ws3.Activate()
$ws3.Range("I30").Select()
